# Rewrites the "Employees" sheet with a new, unstyled data set (as if the
# workbook had been produced fresh by Apache POI instead of hand-formatted
# in Excel): a 4th "Salary" column is added, the roster itself changes, and
# all of the bespoke header styling (bold fill/border, tall header row) is
# stripped back to the workbook default.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original header row (row 1) carries a custom 28pt row height baked in
# from the old formatting. Deleting it outright (rather than just
# overwriting its contents) drops that per-row height so the rebuilt header
# below comes back at the sheet's normal default height.
$ws.Rows(1).Delete()

# Strip any left-over cell styling (the bold/bordered/filled header look,
# and the thin-bordered body cells) from the whole table so every cell goes
# back to the workbook's plain default style.
$ws.Range("A1:D5").ClearFormats()

# Header row
$ws.Range("A1").Value = "Firstname"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "Job_ID"
$ws.Range("D1").Value = "Salary"

# Mary Jones, PO, 200000
$ws.Range("A2").Value = "Mary"
$ws.Range("B2").Value = "Jones"
$ws.Range("C2").Value = "PO"
$ws.Range("D2").Value = 200000

# Vinod Kumar, BA, 110000
$ws.Range("A3").Value = "Vinod"
$ws.Range("B3").Value = "Kumar"
$ws.Range("C3").Value = "BA"
$ws.Range("D3").Value = 110000

# Mansoor Khan, Developer, 135000
$ws.Range("A4").Value = "Mansoor"
$ws.Range("B4").Value = "Khan"
$ws.Range("C4").Value = "Developer"
$ws.Range("D4").Value = 135000

# Linda Smith, SDET, 125000
$ws.Range("A5").Value = "Linda"
$ws.Range("B5").Value = "Smith"
$ws.Range("C5").Value = "SDET"
$ws.Range("D5").Value = 125000

# Reset the view: plain selection back at A1, zoomed in heavily like the
# new file (the workbook was apparently saved while zoomed to 364%).
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 364
